# Applies the "Add files via upload" change: new Apple Parts Stock daily-tracker
# rows on the NOV-2021 sheet (sheet15.xml), pushing the trailing legend rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2021")

# Insert 6 new rows above the existing blank separator row (old row 17),
# which pushes it down to row 23.
$ws.Rows("17:22").Insert()

# The freshly inserted rows pick up the formatting of the row that used to be
# at 17 (the blank separator). Re-apply the look of the data row above
# (row 16) across the 6 new rows so they match the rest of the daily-tracker
# table (one of them - row 17 - stays a blank spacer, like in the source).
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G22").PasteSpecial(-4122)

# Row 18 (14 / 22-Nov-2021)
$ws.Range("A18").Value = 14
$ws.Range("B18").Value = 44522
$ws.Range("C18").Value = "QMVAR 2.0"
$ws.Range("D18").Value = "QMVAR 2.0"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "Completed"
$ws.Range("G18").Value = "Apple Parts Stock (Apple Grid view)"

# Row 19 (15 / 23-Nov-2021)
$ws.Range("A19").Value = 15
$ws.Range("B19").Value = 44523
$ws.Range("C19").Value = "QMVAR 2.0, PENNA"
$ws.Range("D19").Value = "QMVAR 2.0, PENNA"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "Completed"
$ws.Range("G19").Value = "Apple Parts Stock (Add Record), hod role, clinker code move to PROD server"
$ws.Rows("19:19").RowHeight = 28.8

# Row 20 (16 / 24-Nov-2021)
$ws.Range("A20").Value = 16
$ws.Range("B20").Value = 44524
$ws.Range("C20").Value = "QMVAR 2.0"
$ws.Range("D20").Value = "QMVAR 2.0"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Completed"
$ws.Range("G20").Value = "Apple Parts Stock (Design Cheanges)"

# Row 21 (17 / 25-Nov-2021)
$ws.Range("A21").Value = 17
$ws.Range("B21").Value = 44525
$ws.Range("C21").Value = "QMVAR 2.0"
$ws.Range("D21").Value = "QMVAR 2.0"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "Completed"
$ws.Range("G21").Value = "Apple Parts Stock (Sorting, searching)"

# Row 22 (18 / 26-Nov-2021)
$ws.Range("A22").Value = 18
$ws.Range("B22").Value = 44526
$ws.Range("C22").Value = "QMVAR 2.0"
$ws.Range("D22").Value = "QMVAR 2.0"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = "Completed"
$ws.Range("G22").Value = "Apple Parts Stock (Design Cheanges)"

# Keep the view on the freshly entered data, matching the saved selection/scroll.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C18").Select()
